# feat: add 2022-Q4 data
#
# Before: sheet1 "总计" (totals), sheet2 "2022-Q3" (fund holdings for Q3).
# After:  sheet1 "总计" (totals, now with a Q4 row on top of the Q3 row),
#         sheet2 "2022-Q4" (new fund holdings for Q4),
#         sheet3 "2022-Q3" (the original fund-holdings sheet, unchanged).

function Set-TextCell($range, $value) {
    # Force text storage even for numeric-looking strings (fund codes with
    # leading zeros, decimal-looking percentages, etc.) - mirrors Excel's
    # "format cell as Text, then type the value" behaviour.
    $range.NumberFormat = "@"
    $range.Value = $value
}

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsOldQ3 = $wb.Worksheets.Item(2)   # currently "2022-Q3"

# --- 1. Spin the existing "2022-Q3" sheet into the new "2022-Q3" tab -------
# Duplicate the whole sheet object (preserves sheetPr/outline settings,
# formatting, tabSelected, everything) right after itself, then rename the
# original to "2022-Q4" (it will be rewritten with Q4 data below) and the
# duplicate to "2022-Q3". This keeps sheetId/r:id allocation in the 1,2,3
# order matching the sheet tab order, and leaves the duplicate - now named
# "2022-Q3" - as the selected/active sheet, matching the original file.
$wsOldQ3.Copy($null, $wsOldQ3)
$wsOldQ3.Name = "2022-Q4"
$wb.Worksheets.Item(3).Name = "2022-Q3"

# Rename the sheet variable for clarity going forward.
$wsQ4 = $wsOldQ3

# --- 2. Rebuild sheet2 ("2022-Q4") with the new quarter's fund data --------
$wsQ4.Cells.Clear()

$wsTotal.Range("B1").Copy($wsQ4.Range("B1:H1"))

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsTotal.Range("A2").Copy($wsQ4.Range("A2:A3"))

$wsQ4.Range("A2").Value = 0
Set-TextCell $wsQ4.Range("B2") "009649"
$wsQ4.Range("C2").Value = "嘉实精选平衡混合A"
Set-TextCell $wsQ4.Range("D2") "0.47"
Set-TextCell $wsQ4.Range("E2") "68.05"
Set-TextCell $wsQ4.Range("F2") "2.74"
Set-TextCell $wsQ4.Range("G2") "0.0129"
$wsQ4.Range("H2").Value = 9

$wsQ4.Range("A3").Value = 1
Set-TextCell $wsQ4.Range("B3") "009650"
$wsQ4.Range("C3").Value = "嘉实精选平衡混合C"
Set-TextCell $wsQ4.Range("D3") "0.04"
Set-TextCell $wsQ4.Range("E3") "68.05"
Set-TextCell $wsQ4.Range("F3") "2.74"
Set-TextCell $wsQ4.Range("G3") "0.0011"
$wsQ4.Range("H3").Value = 9

# --- 3. Update "总计": insert the Q4 summary row, push Q3 down one row -----
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 0.04

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01
